# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 219 (pushing the existing rows 219-289
# down to 220-290) and populate the new row with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 219:289 down to 220:290, carrying their formatting along.
$ws.Rows("219:219").Insert()

# Populate the newly inserted row 219 with this week's record.
$ws.Range("A219").Value2 = 4
$ws.Range("B219").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C219").Value2 = "Los Lagos"
$ws.Range("D219").Value2 = 44809
$ws.Range("E219").Value2 = 10
$ws.Range("F219").Value2 = 100112032
$ws.Range("G219").Value2 = "Zapallo italiano"
$ws.Range("H219").Value2 = "Sin especificar"
$ws.Range("I219").Value2 = "Primera"
$ws.Range("J219").Value2 = 70
$ws.Range("K219").Value2 = 30000
$ws.Range("L219").Value2 = 30000
$ws.Range("M219").Value2 = 30000
$ws.Range("N219").Value2 = "`$/caja 60 unidades"
$ws.Range("O219").Value2 = "Región de Arica y Parinacota"
$ws.Range("P219").Value2 = 500
$ws.Range("Q219").Value2 = 60
$ws.Range("R219").Value2 = "Hortaliza"
